$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column P (Cle_inf_7_fadiese_g) and column N (Cle_inf_6_mi_g)
# Delete from right to left to keep indices stable.
$ws.Columns.Item(16).Delete()  # P
$ws.Columns.Item(14).Delete()  # N

# Update selection to match the post-edit state
$ws.Range("O15").Select()
